$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 387, shifting existing rows 387:455 down to 388:456
$ws.Rows.Item(387).Insert()

# Fill in the new row 387 with the data (same constant columns as the rest of the data set)
$ws.Cells.Item(387, 1).Value = 4
$ws.Cells.Item(387, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(387, 3).Value = "Los Lagos"
$ws.Cells.Item(387, 4).Value = 45211
$ws.Cells.Item(387, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(387, 5).Value = 10
$ws.Cells.Item(387, 6).Value = 100112044
$ws.Cells.Item(387, 7).Value = "Perejil"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Primera"
$ws.Cells.Item(387, 10).Value = 50
$ws.Cells.Item(387, 11).Value = 6000
$ws.Cells.Item(387, 12).Value = 6000
$ws.Cells.Item(387, 13).Value = 6000
$ws.Cells.Item(387, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(387, 15).Value = "Región Metropolitana"
$ws.Cells.Item(387, 16).Value = 2000
$ws.Cells.Item(387, 17).Value = 3
$ws.Cells.Item(387, 18).Value = "Hortaliza"
